$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pair_matrix")

# The data originally stored in row 29 (Pair 14, Half 2) is reassigned to
# row 30 (Pair 15, Half 1), and vice versa. Columns C (Category), D (Pair)
# and E (Half) stay anchored to their row; every other column (A, B, F,
# and the statistics columns G:Y) swaps between the two rows.

$row29 = @{
    A = 265
    B = "Dog_15"
    F = "Wikipedia"
    G = 3.4313725490196099
    H = 1.74350167692974
    I = 102
    J = 3.9306930693069302
    K = 1.81249786616467
    L = 101
    M = 3.511111111
    N = 1.6871468949999999
    O = 45
    P = 3.4107142860000002
    Q = 1.786584204
    R = 56
    S = 3.7959183670000001
    T = 1.881705696
    U = 49
    V = 4.057692308
    W = 1.7535785429999999
    X = 52
    Y = 0.100396824999999
}

$row30 = @{
    A = 270
    B = "Dog_20"
    F = "Pixabay"
    G = 5.9009900990099
    H = 0.94345058688888905
    I = 101
    J = 4.1941747572815498
    K = 1.63334285167605
    L = 103
    M = 5.6896551720000001
    N = 0.94045335399999996
    O = 58
    P = 6.1860465119999999
    Q = 0.87982178799999999
    R = 43
    S = 4.0816326529999998
    T = 1.5115766420000001
    U = 49
    V = 4.2830188680000001
    W = 1.758248238
    X = 53
    Y = 0.49639133999999902
}

# Write row30's former values into row 29, and row29's former values into row 30.
foreach ($col in $row30.Keys) {
    $ws.Range("$col" + "29").Value = $row30[$col]
}
foreach ($col in $row29.Keys) {
    $ws.Range("$col" + "30").Value = $row29[$col]
}

$ws.Range("E31").Select()
